# Insert a new data row at row 321 (shifts existing rows 321-400 down to 322-401),
# then populate the new row 321 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 321, carrying the rest of the
# table (old rows 321..400) down to 322..401.
$ws.Rows.Item(321).Insert()

# The row that was pushed down (now row 322) kept all of its original values,
# so the new row 321 only needs to be filled in with its own data. Most of the
# descriptive columns repeat the same "Vega Modelo de Temuco" / "Zanahoria"
# values used throughout this sheet.
$ws.Cells.Item(321, 1).Value = 10
$ws.Cells.Item(321, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(321, 3).Value = "La Araucanía"
$ws.Cells.Item(321, 4).Value = 44932
$ws.Cells.Item(321, 5).Value = 9
$ws.Cells.Item(321, 6).Value = 100114013
$ws.Cells.Item(321, 7).Value = "Zanahoria"
$ws.Cells.Item(321, 8).Value = "Sin especificar"
$ws.Cells.Item(321, 9).Value = "Primera"
$ws.Cells.Item(321, 10).Value = 40
$ws.Cells.Item(321, 11).Value = 14000
$ws.Cells.Item(321, 12).Value = 15000
$ws.Cells.Item(321, 13).Value = 14500
$ws.Cells.Item(321, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(321, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(321, 16).Value = 580
$ws.Cells.Item(321, 17).Value = 25
$ws.Cells.Item(321, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(321, 4).NumberFormat = $ws.Cells.Item(322, 4).NumberFormat
